# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.446.56"
$ws.Range("E2").Value = "  -2.48%  "

$ws.Range("D3").Value = "1.789.02"
$ws.Range("E3").Value = "  -2.05%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'231.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("D6").Value = "'0.5870"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.62%  "

$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("D9").Value = "'23.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").Value = "'0.06712"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.40%  "

$ws.Range("D11").Value = "'0.07541"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "

$ws.Range("D12").Value = "1.794.19"
$ws.Range("E12").Value = "  -2.34%  "

$ws.Range("D13").Value = "'4.761"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D14").Value = "'0.6056"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.07%  "

$ws.Range("D15").Value = "2.030.88"
$ws.Range("E15").Value = "  -2.02%  "

$ws.Range("D16").Value = "'75.29"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "'0.000008709"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.59%  "

$ws.Range("D18").Value = "28.418.46"
$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("D19").Value = "'5.400"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.34%  "

$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "'207.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.94%  "

$ws.Range("D22").Value = "'11.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.20%  "

$ws.Range("D23").Value = "'6.767"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.21%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").Value = "'152.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.99%  "

$ws.Range("D26").Value = "'8.055"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.37%  "

$ws.Range("E27").Value = "  -2.44%  "

$ws.Range("D28").Value = "'16.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("D29").Value = "'1.407"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").Value = "'0.06112"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.72%  "

$ws.Range("E31").Value = "  -1.42%  "

$ws.Range("D32").Value = "'3.762"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.87%  "

$ws.Range("D33").Value = "'3.756"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").Value = "'1.673"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.49%  "

$ws.Range("D35").Value = "'1.043"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.41%  "

$ws.Range("D36").Value = "'0.6400"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("D37").Value = "'2.502"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("D38").Value = "'2.696"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.16%  "

$ws.Range("D39").Value = "1.143.11"
$ws.Range("E39").Value = "  -2.21%  "

$ws.Range("D40").Value = "'0.01672"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.48%  "

$ws.Range("D41").Value = "'6.294"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.39%  "

$ws.Range("D42").Value = "'0.8741"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.15%  "

$ws.Range("D43").Value = "'1.006"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("D44").Value = "'100.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("D45").Value = "1.939.22"
$ws.Range("E45").Value = "  -2.01%  "

$ws.Range("D46").Value = "'59.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.37%  "

$ws.Range("D47").Value = "'0.00000000109"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.26%  "

$ws.Range("D48").Value = "'8.389"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.55%  "

$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").Value = "'0.05425"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").Value = "'0.4469"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.74%  "

